$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- zh-cn sheet: Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K) ---
$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e6925de1e760bfa23595e44ccc4666704b871f5/e2e/d8d74326-126c-458b-b46d-c6d636382a8b.md"
$targetMdName = "d8d74326-126c-458b-b46d-c6d636382a8b.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetMdUrl, "", "", $targetMdName)
$wsZhCn.Range("J2").Value = "d8d74326-126c-458b-b46d-c6d636382a8b.b343ca29861dfad8c7e5de96c93573ff1e99ae51.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 01:05:42"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetMdUrl, "", "", $targetMdName)
$wsZhCn.Range("J3").Value = "d8d74326-126c-458b-b46d-c6d636382a8b.b343ca29861dfad8c7e5de96c93573ff1e99ae51.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-20 01:05:42"

# --- de-de sheet: Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K) ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetMdUrl, "", "", $targetMdName)
$wsDeDe.Range("J2").Value = "d8d74326-126c-458b-b46d-c6d636382a8b.b343ca29861dfad8c7e5de96c93573ff1e99ae51.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 01:05:48"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetMdUrl, "", "", $targetMdName)
$wsDeDe.Range("J3").Value = "d8d74326-126c-458b-b46d-c6d636382a8b.b343ca29861dfad8c7e5de96c93573ff1e99ae51.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-20 01:05:48"

# --- Column width adjustments (autofit side-effect of wider content) ---
$wideWidth = 29 + (1/6)
$fullWidth = 39 + (1/6)

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $fullWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fullWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $fullWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fullWidth
